# Adds the "mean vs median" and "variability" Q&A sections (plus a
# preceding page break) to the end of the document, right after the
# paragraph that ends "...a certain campaign or not."

$d = $word.ActiveDocument

# Locate the last paragraph in the main body (the "We could have looked
# at a graph..." paragraph) and build an insertion point right after it.
$lastPara = $d.Paragraphs.Last
$insertAt = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$newXml =
    "<w:p $wNs><w:r><w:br w:type='page'/></w:r></w:p>" +
    "<w:p $wNs><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>" +
    "<w:p $wNs><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Use your data to determine whether the mean or the median better summarizes the data.</w:t></w:r></w:p>" +
    "<w:p $wNs><w:r><w:t>The mean would better represent this data set since there are a lot of counts that are in the thousands and if we were to limit it to just the median, we would be misrepresenting the data as we have records that are way above the median.</w:t></w:r></w:p>" +
    "<w:p $wNs><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Use your data to determine if there is more variability with successful or unsuccessful campaigns. Does this make sense? Why or why not?</w:t></w:r></w:p>" +
    "<w:p $wNs><w:r><w:t xml:space='preserve'>There is more variability in successful campaigns. This does makes sense as to the goal for each campaign is different from one another. Some have low goals while others have very high goals. So, we will see a change in the data provided. Since the goal is not the same across the data set so there would be a fixed pattern in the data since we did not have a same goal or similar goal with the dataset provided. </w:t></w:r></w:p>"

$insertAt.InsertXML($newXml) | Out-Null
